# templateIndividu.xlsx: insert Continent/Country/Ecozone ahead of the
# taxonomy columns (header row gets fully re-ordered) and append four new
# trailing fields (Locality, Number, Collection_Date, Sexe).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing header formatting (bold font, thin border box,
# center/top alignment -> style index 1 in the original file) from A1 and
# stamp it onto the four brand-new trailing header cells before filling
# them in, so Y1:AB1 look like the rest of row 1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("Y1:AB1").PasteSpecial(-4122) | Out-Null

# A1 (SpecimenCode) stays put; B1..X1 are rewritten in the new column
# order, shifting Continent/Country/Ecozone to the front of the taxonomy
# block and pushing everything else right.
$ws.Range("B1").Value  = "Continent"
$ws.Range("C1").Value  = "Country"
$ws.Range("D1").Value  = "Ecozone"
$ws.Range("E1").Value  = "Order"
$ws.Range("F1").Value  = "Suborder"
$ws.Range("G1").Value  = "Tribu"
$ws.Range("H1").Value  = "Family"
$ws.Range("I1").Value  = "Subfamily"
$ws.Range("J1").Value  = "Genus"
$ws.Range("K1").Value  = "Subgenus"
$ws.Range("L1").Value  = "species"
$ws.Range("M1").Value  = "Subspecies"
$ws.Range("N1").Value  = "Num_ID"
$ws.Range("O1").Value  = "Genus_Descriptor"
$ws.Range("P1").Value  = "Species_Descriptor"
$ws.Range("Q1").Value  = "Subgenus_Descriptor"
$ws.Range("R1").Value  = "Subspecies_descriptor"
$ws.Range("S1").Value  = "Genus_Date"
$ws.Range("T1").Value  = "Subgenus_Date"
$ws.Range("U1").Value  = "Species_Date"
$ws.Range("V1").Value  = "Subspecies_Date"
$ws.Range("W1").Value  = "Latitude"
$ws.Range("X1").Value  = "Longitude"

# New trailing columns.
$ws.Range("Y1").Value  = "Locality"
$ws.Range("Z1").Value  = "Number"
$ws.Range("AA1").Value = "Collection_Date"
$ws.Range("AB1").Value = "Sexe"

# Best-effort match of the saved view state (scroll position + selected
# cell) recorded in the sheet's <sheetView>.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
$ws.Range("Z23").Select() | Out-Null
